# Daily update at 8 AM UTC
# Appends the next day's cumulative win counts as a new row, and flips the
# previous "latest" row's date-cell format back to the standard (date+time)
# format now that it's no longer the newest entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column A (the Day column).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

# The previous latest row reverts from the distinctive "today" date-only
# format to the standard date-time format used by all older rows.
$ws.Cells.Item($lastRow, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's data.
$ws.Cells.Item($newRow, 1).Value = 45679
$ws.Cells.Item($newRow, 2).Value = 216
$ws.Cells.Item($newRow, 3).Value = 215
$ws.Cells.Item($newRow, 4).Value = 214

# New row's date cell gets the distinctive "latest entry" date-only format.
$ws.Cells.Item($newRow, 1).NumberFormat = "YYYY-MM-DD"
